$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row from the sheet's used range.
$lastRow = $ws.UsedRange.Rows.Count

# Column K holds "Nombre de centro de carga" (loading terminal name).
# Row 1 is the header; data starts at row 2. Convert each terminal name
# to upper case (accents preserved), matching the diff.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 11)  # Column K = 11
    $val = $cell.Text
    if ($null -ne $val -and $val -ne "") {
        $cell.Value = $val.ToUpper()
    }
}
